# Apply the report-title fix described by the commit "fix the final report title".
#
# The title, the binome's family name and the "Classe" line all changed:
#   - "... Gestion de Livres et Auteurs avec Interface Graphique (Swing)"
#       -> "... Gestion de Clients et Commandes avec Interface Graphique (Swing)"
#   - "Ayoub Majid"  -> "Ayoub Majjid"   (typo fix, doubled "j")
#   - "Classe : "    -> "Classe: " + a tab before "4iir9"

$d = $word.ActiveDocument

# 1) Fix the report title: "Livres et Auteurs" -> "Clients et Commandes"
$d.Content.Find.Execute("Livres et Auteurs", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Clients et Commandes", 2) | Out-Null

# 2) Fix the typo in the family name: "Ayoub Majid" -> "Ayoub Majjid"
$d.Content.Find.Execute("Ayoub Majid", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Ayoub Majjid", 2) | Out-Null

# 3) Tighten "Classe : " to "Classe:" and insert a tab before the class name
$d.Content.Find.Execute("Classe : 4iir9", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Classe: ^t4iir9", 2) | Out-Null
